# Apply the commit's data corrections to the Terminal Status / Cassette
# Balances report. The sheet is a flat table (header row 4, data rows
# 5-27, total row 28); each edit below targets specific cells so that the
# resulting values match the target OOXML exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 6 (ZACATES MARKET) ----
$ws.Range("J6").Value = "07/18/23 17:59"
$ws.Range("N6").Value = 640

# ---- Row 8 becomes LK864765 / SKY LIQUOR ----
$ws.Range("A8").Value = "LK864765"
$ws.Range("C8").Value = "SKY LIQUOR"
$ws.Range("E8").Value = 1560
$ws.Range("I8").Value = ""
$ws.Range("H8").Value = 45130.04188321759
$ws.Range("J8").Value = "07/18/23 18:03"
$ws.Range("K8").Value = "07/18/23 17:46"
$ws.Range("L8").Value = 80
$ws.Range("M8").Value = "$2,420 as of 7/18/2023 8:02:45 AM"
$ws.Range("N8").Value = 1640

# ---- Row 9 becomes L647934 / SB #6 ----
$ws.Range("A9").Value = "L647934"
$ws.Range("C9").Value = "SB #6"
$ws.Range("E9").Value = 1940
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J9").Value = "04/06/23 22:10"
$ws.Range("K9").Value = "04/06/23 22:05"
$ws.Range("L9").Value = 20
$ws.Range("M9").Value = "$1,940 as of 4/6/2023 8:05:45 PM"
$ws.Range("N9").Value = 1960

# ---- Row 10 becomes L704741 / W ADAMS COIN LAUNDRY ----
$ws.Range("A10").Value = "L704741"
$ws.Range("C10").Value = "W ADAMS COIN LAUNDRY"
$ws.Range("E10").Value = 2020
$ws.Range("H10").Value = 45129.04188321759
$ws.Range("J10").Value = "07/18/23 18:59"
$ws.Range("K10").Value = "07/18/23 18:59"
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = "$2,100 as of 7/18/2023 10:08:00 AM"
$ws.Range("N10").Value = 2080

# ---- Row 11 becomes L662336 / SB#4 MONA MARKET ----
$ws.Range("A11").Value = "L662336"
$ws.Range("C11").Value = "SB#4 MONA MARKET"
$ws.Range("E11").Value = 2260
$ws.Range("H11").Value = 45132.04188321759
$ws.Range("J11").Value = "07/18/23 17:18"
$ws.Range("K11").Value = "07/17/23 17:03"
$ws.Range("L11").Value = 100
$ws.Range("M11").Value = "$2,260 as of 7/17/2023 3:03:55 PM"
$ws.Range("N11").Value = 2260

# ---- Row 15 (NICK SHELL SERVICE) ----
$ws.Range("E15").Value = 3540
$ws.Range("J15").Value = "07/18/23 16:42"
$ws.Range("K15").Value = "07/18/23 16:42"
$ws.Range("N15").Value = 3580

# ---- Row 17 (SAFETY MARKET) ----
$ws.Range("E17").Value = 4480
$ws.Range("J17").Value = "07/18/23 17:49"
$ws.Range("K17").Value = "07/18/23 17:28"
$ws.Range("N17").Value = 4480

# ---- Row 18 (DONUT & SANDWICH) ----
$ws.Range("E18").Value = 4960
$ws.Range("J18").Value = "07/18/23 17:43"
$ws.Range("K18").Value = "07/18/23 17:43"
$ws.Range("L18").Value = 20
$ws.Range("N18").Value = 5000

# ---- Row 19 becomes LK236828 / WORLDWIDE AUTOMOTIVE ----
$ws.Range("A19").Value = "LK236828"
$ws.Range("C19").Value = "WORLDWIDE AUTOMOTIVE"
$ws.Range("E19").Value = 5380
$ws.Range("H19").Value = 45151.04188321759
$ws.Range("I19").Value = ""
$ws.Range("J19").Value = "07/18/23 17:06"
$ws.Range("K19").Value = "07/18/23 17:06"
$ws.Range("L19").Value = 80
$ws.Range("M19").Value = "$5,480 as of 7/17/2023 6:02:33 PM"
$ws.Range("N19").Value = 5480

# ---- Row 20 becomes L488595 / N S MART ----
$ws.Range("A20").Value = "L488595"
$ws.Range("C20").Value = "N S MART"
$ws.Range("E20").Value = 5480
$ws.Range("H20").Value = 45285.04188321759
$ws.Range("I20").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J20").Value = "07/16/23 01:58"
$ws.Range("K20").Value = "07/16/23 01:58"
$ws.Range("L20").Value = 60
$ws.Range("M20").Value = "$5,480 as of 7/15/2023 11:58:38 PM"
$ws.Range("N20").Value = 5580

# ---- Row 23 (S.B. 2) ----
$ws.Range("E23").Value = 7400
$ws.Range("J23").Value = "07/18/23 19:24"
$ws.Range("K23").Value = "07/18/23 19:24"
$ws.Range("N23").Value = 7600

# ---- Row 24 (SAMYS PHONE CARDS) ----
$ws.Range("E24").Value = 7960
$ws.Range("J24").Value = "07/18/23 18:58"
$ws.Range("K24").Value = "07/18/23 18:58"
$ws.Range("N24").Value = 8160

# ---- Row 26 (S B DISCOUNT MART) ----
$ws.Range("E26").Value = 11260
$ws.Range("J26").Value = "07/18/23 17:32"
$ws.Range("K26").Value = "07/18/23 17:32"
$ws.Range("N26").Value = 11360

# ---- Row 27 (98 DISCOUNT STORE) ----
$ws.Range("E27").Value = 16220
$ws.Range("J27").Value = "07/18/23 18:26"
$ws.Range("K27").Value = "07/18/23 18:26"
$ws.Range("N27").Value = 16260

# ---- Row 28 (Total Outstanding Cash Balance) ----
$ws.Range("E28").Value = 110600
